# Created function for ensuring the nearest mobile defense is responsible
# for pathing towards the incoming UAS threat.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update effector input values
$ws.Range("A2").Value = 1000
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = 500

$ws.Range("C3").Value = 500

$ws.Range("D4").Value = 0

$ws.Range("D5").Value = 0

$ws.Range("B6").Value = -600
$ws.Range("C6").Value = 500
$ws.Range("D6").Value = 0

# Update selected cell
$ws.Range("D7").Select()

$wb.Save()
